$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E19").Value = "LP L296_EN.pdf"
$wb.Hyperlinks.Add($ws.Range("E19"), "C:\Users\sterl\Downloads\LP L296_EN.pdf", "", "", "C:\Users\sterl\Downloads\LP L296_EN.pdf")
Write-Host "After add: $($ws.Range('E19').Value)"
$ws.Range("E19").Value = "LP L296_EN.pdf"
Write-Host "After reset: $($ws.Range('E19').Value)"
